# Auto-generated edit script: refresh the live crypto price/volume snapshot
# (columns D "Price" and E "Volume(1h)") for data rows 2-51, matching the
# "Updated cryptos list ... with GitHub Actions" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings like "1.030" or "27.513.59" that are NOT valid
# numbers (thousand-grouped with dots, or would lose a trailing zero if Excel
# auto-detected them as numbers). Force Text format before writing so every
# new value round-trips as the exact literal string, just like the original.
$priceCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '27.468.47'
$ws.Range("E2").Value = '  +4.12%  '
$ws.Range("D3").Value = '1.841.99'
$ws.Range("E3").Value = '  +3.72%  '
$ws.Range("D4").Value = '1.026'
$ws.Range("E4").Value = '  +2.57%  '
$ws.Range("D5").Value = '318.99'
$ws.Range("D6").Value = '1.025'
$ws.Range("E6").Value = '  +2.46%  '
$ws.Range("D7").Value = '0.4366'
$ws.Range("E7").Value = '  +3.18%  '
$ws.Range("D8").Value = '0.3730'
$ws.Range("E8").Value = '  +3.69%  '
$ws.Range("D9").Value = '0.07375'
$ws.Range("E9").Value = '  +3.40%  '
$ws.Range("E10").Value = '  +4.62%  '
$ws.Range("D11").Value = '21.41'
$ws.Range("D12").Value = '1.868.13'
$ws.Range("E12").Value = '  +5.23%  '
$ws.Range("D13").Value = '5.476'
$ws.Range("E13").Value = '  +4.51%  '
$ws.Range("D14").Value = '6.696'
$ws.Range("E14").Value = '  +3.87%  '
$ws.Range("D15").Value = '0.07137'
$ws.Range("E15").Value = '  +4.03%  '
$ws.Range("D16").Value = '82.63'
$ws.Range("E16").Value = '  +4.59%  '
$ws.Range("D17").Value = '1.031'
$ws.Range("E17").Value = '  +2.57%  '
$ws.Range("D18").Value = '0.000009005'
$ws.Range("E18").Value = '  +4.33%  '
$ws.Range("D19").Value = '1.026'
$ws.Range("E19").Value = '  +2.49%  '
$ws.Range("D20").Value = '15.41'
$ws.Range("E20").Value = '  +3.53%  '
$ws.Range("D21").Value = '27.486.78'
$ws.Range("E21").Value = '  +4.18%  '
$ws.Range("D22").Value = '5.237'
$ws.Range("E22").Value = '  +3.18%  '
$ws.Range("D23").Value = '11.28'
$ws.Range("E23").Value = '  +2.99%  '
$ws.Range("D24").Value = '2.075.65'
$ws.Range("E24").Value = '  +4.10%  '
$ws.Range("D25").Value = '156.90'
$ws.Range("E25").Value = '  +3.15%  '
$ws.Range("D26").Value = '1.910'
$ws.Range("E26").Value = '  +5.25%  '
$ws.Range("D27").Value = '18.63'
$ws.Range("E27").Value = '  +3.66%  '
$ws.Range("D28").Value = '5.262'
$ws.Range("E28").Value = '  +3.97%  '
$ws.Range("D29").Value = '1.923'
$ws.Range("E29").Value = '  +4.96%  '
$ws.Range("D30").Value = '116.12'
$ws.Range("E30").Value = '  +1.37%  '
$ws.Range("D31").Value = '0.09056'
$ws.Range("E31").Value = '  +2.65%  '
$ws.Range("D32").Value = '1.207'
$ws.Range("E32").Value = '  +7.94%  '
$ws.Range("D33").Value = '0.7622'
$ws.Range("E33").Value = '  +5.06%  '
$ws.Range("D34").Value = '4.488'
$ws.Range("E34").Value = '  +3.92%  '
$ws.Range("D35").Value = '2.871'
$ws.Range("E35").Value = '  +5.10%  '
$ws.Range("D36").Value = '1.028'
$ws.Range("E36").Value = '  +2.79%  '
$ws.Range("D37").Value = '1.148'
$ws.Range("E37").Value = '  +5.07%  '
$ws.Range("D38").Value = '0.01967'
$ws.Range("E38").Value = '  +4.48%  '
$ws.Range("D39").Value = '0.05258'
$ws.Range("E39").Value = '  +2.54%  '
$ws.Range("D40").Value = '0.5174'
$ws.Range("E40").Value = '  +5.47%  '
$ws.Range("D41").Value = '2.777'
$ws.Range("E41").Value = '  +6.80%  '
$ws.Range("D42").Value = '0.1664'
$ws.Range("E42").Value = '  +3.67%  '
$ws.Range("D43").Value = '6.561'
$ws.Range("E43").Value = '  +3.64%  '
$ws.Range("D44").Value = '8.506'
$ws.Range("E44").Value = '  +6.84%  '
$ws.Range("D45").Value = '108.96'
$ws.Range("E45").Value = '  +4.22%  '
$ws.Range("D46").Value = '10.55'
$ws.Range("E46").Value = '  +4.31%  '
$ws.Range("D47").Value = '1.028'
$ws.Range("E47").Value = '  +2.81%  '
$ws.Range("D48").Value = '1.686'
$ws.Range("E48").Value = '  +3.23%  '
$ws.Range("D49").Value = '1.916'
$ws.Range("E49").Value = '  +11.20%  '
$ws.Range("D50").Value = '0.4640'
$ws.Range("E50").Value = '  +4.51%  '
$ws.Range("D51").Value = '0.06314'
$ws.Range("E51").Value = '  +2.30%  '
